# NPC.xlsx / Record_PosList: "add drop item list record"
#
# Row 1 (the header row) on the "Record_PosList" sheet had the generic
# placeholder type name "float" sitting in columns L:O instead of the real
# field names (X / Y / Z / StayTime), while those real names had ended up
# one row too low, in the sample/data row (row 2). This swaps the two rows
# back into place and moves the matching cell comments (which describe the
# item-drop fields) down alongside the values they annotate, and makes the
# "Record_PosList" sheet the active tab/selection instead of "Property".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Property")
$ws2 = $wb.Worksheets.Item("Record_PosList")

# --- 1. Fix the L1:O1 / L2:O2 swap on Record_PosList -----------------------
# Before: row1 = float/float/float/float, row2 = X/Y/Z/StayTime
# After:  row1 = X/Y/Z/StayTime,          row2 = float/float/float/float
$ws2.Range("L1").Value = "X"
$ws2.Range("M1").Value = "Y"
$ws2.Range("N1").Value = "Z"
$ws2.Range("O1").Value = "StayTime"

$ws2.Range("L2").Value = "float"
$ws2.Range("M2").Value = "float"
$ws2.Range("N2").Value = "float"
$ws2.Range("O2").Value = "float"

# --- 2. Move the review comments down onto row 2 with row 1 ----------------
# Capture the text that needs to survive the move before deleting anything.
$commentForL2 = $ws2.Range("M1").Comment.Text()   # "强化等级"
$commentForO2 = $ws2.Range("O1").Comment.Text()   # "镶嵌宝石，逗号分隔"

$ws2.Range("L1").Comment.Delete()   # "物品配置ID" comment is dropped entirely
$ws2.Range("M1").Comment.Delete()
$ws2.Range("N1").Comment.Delete()
$ws2.Range("O1").Comment.Delete()

$ws2.Range("L2").AddComment($commentForL2)
$ws2.Range("M2").AddComment("强化等级")
$ws2.Range("N2").AddComment("强化等级")
$ws2.Range("O2").AddComment($commentForO2)

# --- 3. Update the active sheet/selection -----------------------------------
# "Property" loses tabSelected/its old A40 selection becomes J33 (while it's
# no longer the active sheet); "Record_PosList" becomes the active tab with
# its selection moved from O2 to O8.
[void]$ws1.Activate()
[void]$ws1.Range("J33").Select()

[void]$ws2.Activate()
[void]$ws2.Range("O8").Select()
